$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add U1/V1, copying the bold/bordered header style from T1 ---
$ws.Range("T1").Copy()
$ws.Range("U1:V1").PasteSpecial(-4122)
$ws.Range("U1").Value = "log_elastic_mod_mean"
$ws.Range("V1").Value = "cross_section"

# --- Fill U2:V11 with the new log_elastic_mod_mean / cross_section values ---
$ws.Range("U2").Value = 1.605025141761029
$ws.Range("V2").Value = 21.63104665880347
$ws.Range("U3").Value = 1.944613387595852
$ws.Range("V3").Value = 2.080328294037512
$ws.Range("U4").Value = 1.653724468408931
$ws.Range("V4").Value = 12.32874368763429
$ws.Range("U5").Value = 1.235904591147647
$ws.Range("V5").Value = 204.3424594184416
$ws.Range("U6").Value = 1.047158156233841
$ws.Range("V6").Value = 460.9114030636447
$ws.Range("U7").Value = 1.064938773683436
$ws.Range("V7").Value = 315.4171590574766
$ws.Range("U8").Value = 0.8456531620237618
$ws.Range("V8").Value = 450.5063500201848
$ws.Range("U9").Value = 0.932809100998973
$ws.Range("V9").Value = 875.1097973726238
$ws.Range("U10").Value = 0.538547656797904
$ws.Range("V10").Value = 572.5552611167398
$ws.Range("U11").Value = 1.074221596434249
$ws.Range("V11").Value = 681.1777906100154

# --- New row 12 ("L" series): copy the "series" label style from A11, then fill A:V ---
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("A12").Value = "L"
$ws.Range("B12").Value = 3.5340625
$ws.Range("C12").Value = 0.03035907421343368
$ws.Range("D12").Value = 83.59120132718751
$ws.Range("E12").Value = 12.29472431523172
$ws.Range("F12").Value = 14.70815602602536
$ws.Range("G12").Value = 130
$ws.Range("H12").Value = 70
$ws.Range("I12").Value = 88.66609260977411
$ws.Range("J12").Value = 7.993089968854619
$ws.Range("K12").Value = 37.36461920746228
$ws.Range("L12").Value = 1275.209927873813
$ws.Range("M12").Value = 1443.336022352533
$ws.Range("N12").Value = 20.58771781284712
$ws.Range("O12").Value = 116.4296988845964
$ws.Range("P12").Value = 160.7825106960555
$ws.Range("Q12").Value = 129.4495508456753
$ws.Range("R12").Value = 25.85486850181282
$ws.Range("S12").Value = 2.047709520089304
$ws.Range("T12").Value = 7.920015218587278
$ws.Range("U12").Value = 1.412542333286495
$ws.Range("V12").Value = 9.809307137490865
